$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new note/comment text to E2 (new shared string + new cell)
$ws.Range("E2").Value = "#update ACCEL to take also include the cycles to load within n and catch n! within the CPU"

# Give E2 a thin right border (new border/style entries in styles.xml)
$ws.Range("E2").Borders.Item(10).LineStyle = 1
$ws.Range("E2").Borders.Item(10).Weight = 2

# Nudge the chart (Chart 1) down/right slightly to make room for the new note
$co = $ws.ChartObjects(1)
$co.Left = 348.4375
$co.Top = 14.4
$co.Width = 558.3375
$co.Height = 255.6
